# Apply the "best combination and chains to investigate" update.
#
# Summary of changes:
#  - Outbreak_Locations: add scenario 5 and 6 case locations (rows 24-29)
#  - Store_Locations: add scenario 5 and 6 store locations for Chain 1 and
#    the new "Chain 2" (rows 13-24)
#  - Population: add "min"/"max" columns (D/E), fill in scenario 4's blank
#    population_per_cell cell with a highlight, and add scenario 5 and 6
#    population rows (rows 6-7)
#  - Leaves the Population sheet as the active/selected sheet, matching the
#    author's last working view.

$wb = $excel.ActiveWorkbook

$wsOutbreak = $wb.Worksheets.Item("Outbreak_Locations")
$wsStores   = $wb.Worksheets.Item("Store_Locations")
$wsPop      = $wb.Worksheets.Item("Population")

# ---------------------------------------------------------------------
# Outbreak_Locations: new case rows for scenario 5 and scenario 6
# ---------------------------------------------------------------------
$outbreakRows = @(
  @(5, 150, 850),
  @(5, 250, 750),
  @(5, 350, 850),
  @(6, 550, 550),
  @(6, 250, 750),
  @(6, 550, 850)
)

$r = 24
foreach ($row in $outbreakRows) {
  $wsOutbreak.Cells.Item($r, 1).Value = $row[0]
  $wsOutbreak.Cells.Item($r, 2).Value = $row[1]
  $wsOutbreak.Cells.Item($r, 3).Value = $row[2]
  $r++
}

# ---------------------------------------------------------------------
# Store_Locations: new store rows for scenario 5 and scenario 6,
# including the new "Chain 2"
# ---------------------------------------------------------------------
$storeRows = @(
  @(5, 115, 825, "Chain 1"),
  @(5, 230, 735, "Chain 1"),
  @(5, 345, 855, "Chain 1"),
  @(5, 120, 820, "Chain 2"),
  @(5, 235, 730, "Chain 2"),
  @(5, 340, 850, "Chain 2"),
  @(6, 230, 735, "Chain 1"),
  @(6, 223, 245, "Chain 1"),
  @(6, 523, 822, "Chain 1"),
  @(6, 546, 555, "Chain 1"),
  @(6, 233, 735, "Chain 2"),
  @(6, 524, 829, "Chain 2")
)

$r = 13
foreach ($row in $storeRows) {
  $wsStores.Cells.Item($r, 1).Value = $row[0]
  $wsStores.Cells.Item($r, 2).Value = $row[1]
  $wsStores.Cells.Item($r, 3).Value = $row[2]
  $wsStores.Cells.Item($r, 4).Value = $row[3]
  $r++
}

# ---------------------------------------------------------------------
# Population: new "min" / "max" columns, highlight the blank
# population_per_cell cell for scenario 4, and add scenario 5/6 rows
# ---------------------------------------------------------------------
$wsPop.Range("D1").Value = "min"
$wsPop.Range("E1").Value = "max"

# Scenario 4 uses "random" population, so population_per_cell (C5) stays
# blank but gets highlighted, and the new min/max bounds are filled in.
$wsPop.Range("C5").Interior.Color = 15132391
$wsPop.Range("D5").Value = 1
$wsPop.Range("E5").Value = 100

$wsPop.Range("A6").Value = 5
$wsPop.Range("B6").Value = "uniform"
$wsPop.Range("C6").Value = 5

$wsPop.Range("A7").Value = 6
$wsPop.Range("B7").Value = "uniform"
$wsPop.Range("C7").Value = 5

# ---------------------------------------------------------------------
# Selections / active sheet: the author ended up on the Population
# sheet with E5 selected; leave the other two sheets' last selections
# where the author left them.
# ---------------------------------------------------------------------
$wsOutbreak.Range("C29").Select()
$wsStores.Range("F7").Select()

$wsPop.Activate()
$wsPop.Range("E5").Select()
